$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backup codes for rows 2-4 (3 codes)
$ws.Range("A2").Value = "QS3W554CY3ZX"
$ws.Range("A3").Value = "Q8YK0WWBPD6H"
$ws.Range("A4").Value = "YYZGQ1P8K4EG"

# Clear old rows 14-17
$ws.Range("A14:A17").ClearContents()

# New backup codes for rows 11-16 (6 codes)
$ws.Range("A11").Value = "3Z6ADAYX8TXX"
$ws.Range("A12").Value = "96R6XXZ5H6HD"
$ws.Range("A13").Value = "51XZBTACEQGM"
$ws.Range("A14").Value = "T91KP1V5CTSN"
$ws.Range("A15").Value = "YW0AXXBBWGJP"
$ws.Range("A16").Value = "KBN7HS57G4H3"

$ws.Range("A4").Select()
